$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove all existing hyperlinks first (will rebuild after values are set).
$ws.Hyperlinks.Delete()

# Write A:E, G, H for rows 2-19 (row 20 does not exist in the target).
$ws.Cells.Item(2, 1).Value = '2025-10-11 01:13:47'
$ws.Cells.Item(2, 2).Value = 'マッチングアプリのAIレコメンド構築'
$ws.Cells.Item(2, 3).Value = 'システム開発'
$ws.Cells.Item(2, 4).Value = '100,000 円 ~ 200,000 円 / 固定'
$ws.Cells.Item(2, 5).Value = '期限情報なし'
$ws.Cells.Item(2, 6).Value = 'https://www.lancers.jp/work/detail/5410515'
$ws.Cells.Item(2, 7).Value = 338
$ws.Cells.Item(2, 8).Value = '🔥AI,Ai ◇アプリ'

$ws.Cells.Item(3, 1).Value = '2025-10-11 01:13:47'
$ws.Cells.Item(3, 2).Value = '急募 Zoom/Meet×TLDV×ChatGPT×Notion×Slack 議事録ワークフロー構築依頼'
$ws.Cells.Item(3, 3).Value = 'システム開発'
$ws.Cells.Item(3, 4).Value = '50,000 円 ~ 100,000 円 / 固定'
$ws.Cells.Item(3, 5).Value = '期限情報なし'
$ws.Cells.Item(3, 6).Value = 'https://www.lancers.jp/work/detail/5410688'
$ws.Cells.Item(3, 7).Value = 323
$ws.Cells.Item(3, 8).Value = '🔥GPT,ChatGPT'

$ws.Cells.Item(4, 1).Value = '2025-10-11 01:13:47'
$ws.Cells.Item(4, 2).Value = '【せどり×ツール製作】APIを使用したせどりツールを製作できるエンジニアさんを募集します♪'
$ws.Cells.Item(4, 3).Value = 'システム開発'
$ws.Cells.Item(4, 4).Value = '20,000 円 ~ 50,000 円 / 固定'
$ws.Cells.Item(4, 5).Value = '期限情報なし'
$ws.Cells.Item(4, 6).Value = 'https://www.lancers.jp/work/detail/5217096'
$ws.Cells.Item(4, 7).Value = 243
$ws.Cells.Item(4, 8).Value = '🔥API ◆ツール'

$ws.Cells.Item(5, 1).Value = '2025-10-11 01:13:47'
$ws.Cells.Item(5, 2).Value = '【急募】配送状況を自動取得するAPI開発者募集!'
$ws.Cells.Item(5, 3).Value = 'システム開発'
$ws.Cells.Item(5, 4).Value = '20,000 円 ~ 50,000 円 / 固定'
$ws.Cells.Item(5, 5).Value = '期限情報なし'
$ws.Cells.Item(5, 6).Value = 'https://www.lancers.jp/work/detail/5411268'
$ws.Cells.Item(5, 7).Value = 238
$ws.Cells.Item(5, 8).Value = '🔥API ◆開発'

$ws.Cells.Item(6, 1).Value = '2025-10-11 01:13:47'
$ws.Cells.Item(6, 2).Value = 'Google Apps ScriptとAPIを使用したサイボウズOfficeの連携システム改修'
$ws.Cells.Item(6, 3).Value = 'システム開発'
$ws.Cells.Item(6, 4).Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Cells.Item(6, 5).Value = '期限情報なし'
$ws.Cells.Item(6, 6).Value = 'https://www.lancers.jp/work/detail/5410872'
$ws.Cells.Item(6, 7).Value = 210
$ws.Cells.Item(6, 8).Value = '🔥API'

$ws.Cells.Item(7, 1).Value = '2025-10-11 01:13:47'
$ws.Cells.Item(7, 2).Value = 'LINE WORKSで「URLを外部ブラウザで開く」設定を API経由でON にしてください。'
$ws.Cells.Item(7, 3).Value = 'システム開発'
$ws.Cells.Item(7, 4).Value = '5,000 円 ~ 10,000 円 / 固定'
$ws.Cells.Item(7, 5).Value = '期限情報なし'
$ws.Cells.Item(7, 6).Value = 'https://www.lancers.jp/work/detail/5410829'
$ws.Cells.Item(7, 7).Value = 180
$ws.Cells.Item(7, 8).Value = '🔥API'

$ws.Cells.Item(8, 1).Value = '2025-10-11 01:13:47'
$ws.Cells.Item(8, 2).Value = '【相談希望】在庫管理・出品補助ツールの開発に関するZoom面談依頼'
$ws.Cells.Item(8, 3).Value = 'システム開発'
$ws.Cells.Item(8, 4).Value = '20,000 円 ~ 50,000 円 / 固定'
$ws.Cells.Item(8, 5).Value = '期限情報なし'
$ws.Cells.Item(8, 6).Value = 'https://www.lancers.jp/work/detail/5398112'
$ws.Cells.Item(8, 7).Value = 158
$ws.Cells.Item(8, 8).Value = '◆ツール,開発 ◇管理'

$ws.Cells.Item(9, 1).Value = '2025-10-11 01:13:47'
$ws.Cells.Item(9, 2).Value = '【システム開発】顧客予約サインシステムの構築依頼'
$ws.Cells.Item(9, 3).Value = 'システム開発'
$ws.Cells.Item(9, 4).Value = '20,000 円 ~ 50,000 円 / 固定'
$ws.Cells.Item(9, 5).Value = '期限情報なし'
$ws.Cells.Item(9, 6).Value = 'https://www.lancers.jp/work/detail/5410801'
$ws.Cells.Item(9, 7).Value = 113
$ws.Cells.Item(9, 8).Value = '◆開発,システム開発'

$ws.Cells.Item(10, 1).Value = '2025-10-11 01:13:47'
$ws.Cells.Item(10, 2).Value = '【新規教育プラットフォーム開発】ノーコード・ローコードで構築できる学習アプリ開発パートナー募集!'
$ws.Cells.Item(10, 3).Value = 'システム開発'
$ws.Cells.Item(10, 4).Value = '100,000 円 ~ 200,000 円 / 固定'
$ws.Cells.Item(10, 5).Value = '期限情報なし'
$ws.Cells.Item(10, 6).Value = 'https://www.lancers.jp/work/detail/5410616'
$ws.Cells.Item(10, 7).Value = 93
$ws.Cells.Item(10, 8).Value = '◆開発 ◇アプリ'

$ws.Cells.Item(11, 1).Value = '2025-10-11 01:13:47'
$ws.Cells.Item(11, 2).Value = '【急募】WEBシステムのデータ解析レポート出力開発'
$ws.Cells.Item(11, 3).Value = 'システム開発'
$ws.Cells.Item(11, 4).Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Cells.Item(11, 5).Value = '期限情報なし'
$ws.Cells.Item(11, 6).Value = 'https://www.lancers.jp/work/detail/5410793'
$ws.Cells.Item(11, 7).Value = 90
$ws.Cells.Item(11, 8).Value = '◆開発'

$ws.Cells.Item(12, 1).Value = '2025-10-11 01:13:47'
$ws.Cells.Item(12, 2).Value = '大手クレジットカード企業向け、Google Cloudを利用したアジャイル開発共通基盤案件'
$ws.Cells.Item(12, 3).Value = 'システム開発'
$ws.Cells.Item(12, 4).Value = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Cells.Item(12, 5).Value = '期限情報なし'
$ws.Cells.Item(12, 6).Value = 'https://www.lancers.jp/work/detail/5410520'
$ws.Cells.Item(12, 7).Value = 75
$ws.Cells.Item(12, 8).Value = '◆開発'

$ws.Cells.Item(13, 1).Value = '2025-10-11 01:13:47'
$ws.Cells.Item(13, 2).Value = '大手クレジットカード企業向け、Google Cloudを利用したアジャイル開発共通基盤案件_ワーカー'
$ws.Cells.Item(13, 3).Value = 'システム開発'
$ws.Cells.Item(13, 4).Value = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Cells.Item(13, 5).Value = '期限情報なし'
$ws.Cells.Item(13, 6).Value = 'https://www.lancers.jp/work/detail/5410523'
$ws.Cells.Item(13, 7).Value = 75
$ws.Cells.Item(13, 8).Value = '◆開発'

$ws.Cells.Item(14, 1).Value = '2025-10-11 01:13:47'
$ws.Cells.Item(14, 2).Value = 'スプレッドシートをもとにした顧客・売上管理アプリのグライド化(Glide/無料版)'
$ws.Cells.Item(14, 3).Value = 'システム開発'
$ws.Cells.Item(14, 4).Value = '5,000 円 ~ 10,000 円 / 固定'
$ws.Cells.Item(14, 5).Value = '期限情報なし'
$ws.Cells.Item(14, 6).Value = 'https://www.lancers.jp/work/detail/5411304'
$ws.Cells.Item(14, 7).Value = 55
$ws.Cells.Item(14, 8).Value = '◇アプリ'

$ws.Cells.Item(15, 1).Value = '2025-10-11 01:13:47'
$ws.Cells.Item(15, 2).Value = '【急募】salamに関するウェブサイト制作の依頼'
$ws.Cells.Item(15, 3).Value = 'システム開発'
$ws.Cells.Item(15, 4).Value = '100,000 円 ~ 200,000 円 / 固定'
$ws.Cells.Item(15, 5).Value = '期限情報なし'
$ws.Cells.Item(15, 6).Value = 'https://www.lancers.jp/work/detail/5411046'
$ws.Cells.Item(15, 7).Value = 38
$ws.Cells.Item(15, 8).Value = '◇サイト'

$ws.Cells.Item(16, 1).Value = '2025-10-11 01:13:47'
$ws.Cells.Item(16, 2).Value = '【急募】16タイプ診断コンテンツのLP制作'
$ws.Cells.Item(16, 3).Value = 'システム開発'
$ws.Cells.Item(16, 4).Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Cells.Item(16, 5).Value = '期限情報なし'
$ws.Cells.Item(16, 6).Value = 'https://www.lancers.jp/work/detail/5408735'
$ws.Cells.Item(16, 7).Value = 25
$ws.Cells.Item(16, 8).ClearContents()

$ws.Cells.Item(17, 1).Value = '2025-10-11 01:13:47'
$ws.Cells.Item(17, 2).Value = '〖リモート可〗Delphiエンジニア募集'
$ws.Cells.Item(17, 3).Value = 'システム開発'
$ws.Cells.Item(17, 4).Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Cells.Item(17, 5).Value = '期限情報なし'
$ws.Cells.Item(17, 6).Value = 'https://www.lancers.jp/work/detail/5341051'
$ws.Cells.Item(17, 7).Value = 25
$ws.Cells.Item(17, 8).ClearContents()

$ws.Cells.Item(18, 1).Value = '2025-10-11 01:13:47'
$ws.Cells.Item(18, 2).Value = '初回 【継続あり】Microsoft PL-300/400/600 資格試験向け問題集作成'
$ws.Cells.Item(18, 3).Value = 'システム開発'
$ws.Cells.Item(18, 4).Value = '50,000 円 ~ 100,000 円 / 固定'
$ws.Cells.Item(18, 5).Value = '期限情報なし'
$ws.Cells.Item(18, 6).Value = 'https://www.lancers.jp/work/detail/5411149'
$ws.Cells.Item(18, 7).Value = 18
$ws.Cells.Item(18, 8).ClearContents()

$ws.Cells.Item(19, 1).Value = '2025-10-11 01:13:47'
$ws.Cells.Item(19, 2).Value = '急募 【緊急】selenium(ruby)でのX自動ログインの実装'
$ws.Cells.Item(19, 3).Value = 'システム開発'
$ws.Cells.Item(19, 4).Value = '50,000 円 ~ 100,000 円 / 固定'
$ws.Cells.Item(19, 5).Value = '期限情報なし'
$ws.Cells.Item(19, 6).Value = 'https://www.lancers.jp/work/detail/5411088'
$ws.Cells.Item(19, 7).Value = 18
$ws.Cells.Item(19, 8).ClearContents()

# Restore the Hyperlink cell style on column F (Hyperlinks.Delete() above
# bumped the style index on cells that previously carried it away from
# the shared 'Hyperlink' cellXfs entry).
for ($r = 2; $r -le 19; $r++) {
    $ws.Cells.Item($r, 6).Style = "Hyperlink"
}

# Re-create the hyperlinks themselves, in row order, pointing at the URL text.
$ws.Hyperlinks.Add($ws.Range('F2'), 'https://www.lancers.jp/work/detail/5410515')
$ws.Hyperlinks.Add($ws.Range('F3'), 'https://www.lancers.jp/work/detail/5410688')
$ws.Hyperlinks.Add($ws.Range('F4'), 'https://www.lancers.jp/work/detail/5217096')
$ws.Hyperlinks.Add($ws.Range('F5'), 'https://www.lancers.jp/work/detail/5411268')
$ws.Hyperlinks.Add($ws.Range('F6'), 'https://www.lancers.jp/work/detail/5410872')
$ws.Hyperlinks.Add($ws.Range('F7'), 'https://www.lancers.jp/work/detail/5410829')
$ws.Hyperlinks.Add($ws.Range('F8'), 'https://www.lancers.jp/work/detail/5398112')
$ws.Hyperlinks.Add($ws.Range('F9'), 'https://www.lancers.jp/work/detail/5410801')
$ws.Hyperlinks.Add($ws.Range('F10'), 'https://www.lancers.jp/work/detail/5410616')
$ws.Hyperlinks.Add($ws.Range('F11'), 'https://www.lancers.jp/work/detail/5410793')
$ws.Hyperlinks.Add($ws.Range('F12'), 'https://www.lancers.jp/work/detail/5410520')
$ws.Hyperlinks.Add($ws.Range('F13'), 'https://www.lancers.jp/work/detail/5410523')
$ws.Hyperlinks.Add($ws.Range('F14'), 'https://www.lancers.jp/work/detail/5411304')
$ws.Hyperlinks.Add($ws.Range('F15'), 'https://www.lancers.jp/work/detail/5411046')
$ws.Hyperlinks.Add($ws.Range('F16'), 'https://www.lancers.jp/work/detail/5408735')
$ws.Hyperlinks.Add($ws.Range('F17'), 'https://www.lancers.jp/work/detail/5341051')
$ws.Hyperlinks.Add($ws.Range('F18'), 'https://www.lancers.jp/work/detail/5411149')
$ws.Hyperlinks.Add($ws.Range('F19'), 'https://www.lancers.jp/work/detail/5411088')
